$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the extra header columns (B1:G1) ------------------------------
# Only "Cedula" (A1) is kept; the rest of the header row (Nombre Completo,
# Fecha Nacimiento, Nacionalidad, Padre, Madre, Edad) is dropped entirely so
# only a single shared string remains in the workbook.
$ws.Range("B1:G1").Clear()

# A1 also loses the thin-border style that used to box in the header row.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = "Cedula"

# --- Update / add the "Cedula" data rows -----------------------------------
$ws.Range("A2").Value = 117100591
$ws.Range("A3").Value = 117100592
$ws.Range("A4").Value = 117100593
$ws.Range("A5").Value = 117100594

# These cells keep/gain the underline placeholder style used elsewhere in
# the sheet (maps to the same cellXfs entry as B2/D2/D5/B16/C16).
$ws.Range("C4").Font.Underline = 2
$ws.Range("A10").Font.Underline = 2

# --- Selection / view state --------------------------------------------
[void]$ws.Range("B7").Select()

# --- Page setup ------------------------------------------------------------
$ws.PageSetup.Orientation = 1
